$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F ("想去人数" / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4443
$ws1.Range("F3").Value  = 845
$ws1.Range("F5").Value  = 138
$ws1.Range("F6").Value  = 153
$ws1.Range("F7").Value  = 34
$ws1.Range("F8").Value  = 19
$ws1.Range("F9").Value  = 141
$ws1.Range("F10").Value = 616
$ws1.Range("F12").Value = 189
$ws1.Range("F13").Value = 1220
$ws1.Range("F15").Value = 2837
$ws1.Range("F16").Value = 442
$ws1.Range("F17").Value = 528

# Sheet "全部类型" (All types) - same column F updates, shifted by one extra row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4443
$ws4.Range("F3").Value  = 845
$ws4.Range("F5").Value  = 138
$ws4.Range("F6").Value  = 153
$ws4.Range("F7").Value  = 34
$ws4.Range("F8").Value  = 19
$ws4.Range("F9").Value  = 141
$ws4.Range("F10").Value = 616
$ws4.Range("F13").Value = 189
$ws4.Range("F14").Value = 1220
$ws4.Range("F16").Value = 2837
$ws4.Range("F17").Value = 442
$ws4.Range("F18").Value = 528
